$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quote row for 2025-11-24 (Excel serial date 45985).
# Match the number format already used by the preceding date cells (column A).
$ws.Range("A80").NumberFormat = $ws.Range("A79").NumberFormat

$ws.Range("A80").Value = 45985
$ws.Range("B80").Value = "21,2604"
$ws.Range("C80").Value = "15,5723"
$ws.Range("D80").Value = "15,1879"
$ws.Range("E80").Value = "15,1879"
